$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 209; this shifts existing rows 209..318 down to 210..319
$ws.Rows(209).Insert()

# Populate the newly inserted row 209 with the new data record.
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F Categoria ID,
# G Categoria, H Variedad, I Calidad, J Volumen, K Precio minimo, L Precio maximo,
# M Precio promedio ponderado, N Unidad de comercializacion, O Origen,
# P Precio $/Kg, Q Kg o Unidades, R Clasificacion
$ws.Cells.Item(209, 1).Value = 8
$ws.Cells.Item(209, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(209, 3).Value = "Coquimbo"
$ws.Cells.Item(209, 4).Value = 44917
$ws.Cells.Item(209, 5).Value = 4
$ws.Cells.Item(209, 6).Value = 100112031
$ws.Cells.Item(209, 7).Value = "Poroto verde"
$ws.Cells.Item(209, 8).Value = "Magnum"
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 520
$ws.Cells.Item(209, 11).Value = 23000
$ws.Cells.Item(209, 12).Value = 24000
$ws.Cells.Item(209, 13).Value = 23500
$ws.Cells.Item(209, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(209, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(209, 16).Value = 940
$ws.Cells.Item(209, 17).Value = 25
$ws.Cells.Item(209, 18).Value = "Hortaliza"
